$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.043437242507935
$ws.Range("B1").Value = 2.417041301727295
$ws.Range("C1").Value = 5.225387096405029
$ws.Range("D1").Value = 2.303689241409302
$ws.Range("E1").Value = 1.323248028755188
